$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.325.32"
$ws.Range("E2").Value = "  -3.63%  "

$ws.Range("D3").Value = "3.124.24"
$ws.Range("E3").Value = "  -4.57%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "560.00"
$ws.Range("E5").Value = "  -4.37%  "

$ws.Range("D6").Value = "161.31"
$ws.Range("E6").Value = "  -9.15%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -8.70%  "

$ws.Range("D9").Value = "3.117.99"
$ws.Range("E9").Value = "  -4.81%  "

$ws.Range("E10").Value = "  -1.98%  "

$ws.Range("E11").Value = "  -7.33%  "

$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -5.01%  "

$ws.Range("D13").Value = "3.663.73"
$ws.Range("E13").Value = "  -4.79%  "

$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").Value = "63.359.77"
$ws.Range("E15").Value = "  -3.77%  "

$ws.Range("D16").Value = "24.73"
$ws.Range("E16").Value = "  -5.58%  "

$ws.Range("D17").Value = "3.127.84"
$ws.Range("E17").Value = "  -5.62%  "

$ws.Range("E18").Value = "  -6.05%  "

$ws.Range("D19").Value = "398.67"
$ws.Range("E19").Value = "  -5.17%  "

$ws.Range("E20").Value = "  -4.62%  "

$ws.Range("D21").Value = "12.41"
$ws.Range("E21").Value = "  -4.18%  "

$ws.Range("E22").Value = "  -2.83%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").Value = "67.32"
$ws.Range("E24").Value = "  -5.17%  "

$ws.Range("E25").Value = "  -2.98%  "

$ws.Range("E26").Value = "  -5.34%  "

$ws.Range("D27").Value = "0.0₃0999"
$ws.Range("E27").Value = "  -11.05%  "

$ws.Range("E28").Value = "  -7.24%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  -6.76%  "

$ws.Range("D32").Value = "20.80"
$ws.Range("E32").Value = "  -5.78%  "

$ws.Range("E33").Value = "  -4.55%  "

$ws.Range("D34").Value = "4.74"
$ws.Range("E34").Value = "  -6.90%  "

$ws.Range("E35").Value = "  -6.57%  "

$ws.Range("D36").Value = "152.33"
$ws.Range("E36").Value = "  -3.64%  "

$ws.Range("E37").Value = "  -7.59%  "

$ws.Range("D38").Value = "2.741.71"
$ws.Range("E38").Value = "  -3.52%  "

$ws.Range("E39").Value = "  -7.78%  "

$ws.Range("D40").Value = "4.04"
$ws.Range("E40").Value = "  -6.20%  "

$ws.Range("D41").Value = "23.20"
$ws.Range("E41").Value = "  -10.77%  "

$ws.Range("D42").Value = "38.55"
$ws.Range("E42").Value = "  -2.43%  "

$ws.Range("E43").Value = "  -6.91%  "

$ws.Range("D44").Value = "0.0614"
$ws.Range("E44").Value = "  -3.11%  "

$ws.Range("D45").Value = "5.35"
$ws.Range("E45").Value = "  -8.23%  "

$ws.Range("E46").Value = "  -4.07%  "

$ws.Range("D47").Value = "20.68"
$ws.Range("E47").Value = "  -8.69%  "

$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").Value = "280.20"
$ws.Range("E49").Value = "  -9.64%  "

$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  -4.50%  "

$ws.Range("D51").Value = "10.46"
$ws.Range("E51").Value = "  +0.85%  "
